$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -7
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = -2
$ws.Range("F8").Value = -4
$ws.Range("F12").Value = 0
